$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared-string runs) ---
$ws.Range("A8").Value = "Volume 32   Number  32"
$ws.Range("C9").Value = "Report Covering the Week  8/4/2025  Through  8/10/2025"

# --- Donor cells for text ("n/a"-style) values, reused via Copy to preserve exact cell style ---
# C33/D33/G33 hold the text "0" (shared string 20) with style 13
# E33/H33 hold the text "***.*" (shared string 21) with style 13

# --- Data cell updates, rows 14-31 ---
$ws.Range("L14").Value = -33.333333333333
$ws.Range("C33").Copy($ws.Range("D15"))
$ws.Range("E33").Copy($ws.Range("E15"))
$ws.Range("F15").Value = 2
$ws.Range("H15").Value = 100
$ws.Range("I15").Value = 13
$ws.Range("K15").Value = 116.666666666667
$ws.Range("L15").Value = 333.333333333333
$ws.Range("M15").Value = 333.333333333333
$ws.Range("N15").Value = 160
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 9
$ws.Range("G16").Value = 25
$ws.Range("H16").Value = -64
$ws.Range("I16").Value = 86
$ws.Range("J16").Value = 93
$ws.Range("K16").Value = -7.526881720430
$ws.Range("L16").Value = 8.860759493670
$ws.Range("M16").Value = 50.877192982456
$ws.Range("N16").Value = -82.555780933062
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 500
$ws.Range("F17").Value = 26
$ws.Range("G17").Value = 19
$ws.Range("H17").Value = 36.842105263157
$ws.Range("I17").Value = 153
$ws.Range("J17").Value = 130
$ws.Range("K17").Value = 17.692307692307
$ws.Range("L17").Value = 24.390243902439
$ws.Range("M17").Value = 80
$ws.Range("N17").Value = -16.393442622950
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -33.333333333333
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = 62.5
$ws.Range("I18").Value = 81
$ws.Range("J18").Value = 86
$ws.Range("K18").Value = -5.813953488372
$ws.Range("L18").Value = -13.829787234042
$ws.Range("M18").Value = 5.194805194805
$ws.Range("N18").Value = -82.197802197802
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = -36.363636363636
$ws.Range("F19").Value = 39
$ws.Range("G19").Value = 58
$ws.Range("H19").Value = -32.758620689655
$ws.Range("I19").Value = 410
$ws.Range("J19").Value = 443
$ws.Range("K19").Value = -7.449209932279
$ws.Range("L19").Value = -2.843601895734
$ws.Range("M19").Value = 32.258064516129
$ws.Range("N19").Value = -56.977964323189
$ws.Range("C20").Value = 1
$ws.Range("I20").Value = 11
$ws.Range("K20").Value = -31.25
$ws.Range("L20").Value = -62.068965517241
$ws.Range("M20").Value = 22.222222222222
$ws.Range("N20").Value = -93.989071038251
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 21
$ws.Range("E21").Value = -4.761904761904
$ws.Range("F21").Value = 90
$ws.Range("G21").Value = 112
$ws.Range("H21").Value = -19.642857142857
$ws.Range("I21").Value = 756
$ws.Range("J21").Value = 774
$ws.Range("K21").Value = -2.325581395348
$ws.Range("L21").Value = 0.398406374501
$ws.Range("M21").Value = 39.741219963031
$ws.Range("N21").Value = -66.783831282952
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = 100
$ws.Range("F22").Value = 5
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 150
$ws.Range("I22").Value = 39
$ws.Range("J22").Value = 32
$ws.Range("K22").Value = 21.875
$ws.Range("L22").Value = 85.714285714285
$ws.Range("M22").Value = 30
$ws.Range("C33").Copy($ws.Range("C23"))
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = -100
$ws.Range("F23").Value = 3
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 37
$ws.Range("J23").Value = 30
$ws.Range("K23").Value = 23.333333333333
$ws.Range("L23").Value = 54.166666666666
$ws.Range("M23").Value = 94.736842105263
$ws.Range("C24").Value = 26
$ws.Range("D24").Value = 28
$ws.Range("E24").Value = -7.142857142857
$ws.Range("F24").Value = 97
$ws.Range("G24").Value = 115
$ws.Range("H24").Value = -15.652173913043
$ws.Range("I24").Value = 819
$ws.Range("J24").Value = 803
$ws.Range("K24").Value = 1.992528019925
$ws.Range("L24").Value = 33.823529411764
$ws.Range("M24").Value = 20.087976539589
$ws.Range("C25").Value = 15
$ws.Range("D25").Value = 30
$ws.Range("E25").Value = -50
$ws.Range("F25").Value = 58
$ws.Range("G25").Value = 112
$ws.Range("H25").Value = -48.214285714285
$ws.Range("I25").Value = 612
$ws.Range("J25").Value = 686
$ws.Range("K25").Value = -10.787172011661
$ws.Range("L25").Value = 56.122448979591
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 20
$ws.Range("F26").Value = 20
$ws.Range("H26").Value = -42.857142857142
$ws.Range("I26").Value = 260
$ws.Range("J26").Value = 277
$ws.Range("K26").Value = -6.137184115523
$ws.Range("L26").Value = 17.117117117117
$ws.Range("M26").Value = 50.289017341040
$ws.Range("C27").Value = 2
$ws.Range("C33").Copy($ws.Range("D27"))
$ws.Range("E33").Copy($ws.Range("E27"))
$ws.Range("F27").Value = 3
$ws.Range("H27").Value = 200
$ws.Range("I27").Value = 16
$ws.Range("K27").Value = 23.076923076923
$ws.Range("L27").Value = 100
$ws.Range("C33").Copy($ws.Range("C28"))
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 3
$ws.Range("G28").Value = 7
$ws.Range("H28").Value = -57.142857142857
$ws.Range("I28").Value = 41
$ws.Range("J28").Value = 45
$ws.Range("K28").Value = -8.888888888888
$ws.Range("L28").Value = 7.894736842105
$ws.Range("C33").Copy($ws.Range("D29"))
$ws.Range("E33").Copy($ws.Range("E29"))
$ws.Range("C33").Copy($ws.Range("D30"))
$ws.Range("E33").Copy($ws.Range("E30"))
$ws.Range("C33").Copy($ws.Range("D31"))
$ws.Range("E33").Copy($ws.Range("E31"))
$ws.Range("L31").Value = -40
